{"js": "const body = context.document.body;\n\n// 1. Expand \"...built from 2011 to 2016 through the United States\" into\n//    \"...built from 2011 to 2016, throughout the United States,\" (adds\n//    \", \" before \"through\", turns it into \"throughout\", and appends a\n//    trailing comma) while leaving the rest of the paragraph untouched.\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\nlet paragraphs = body.paragraphs.items;\nconst introParagraph = paragraphs[1];\n\nconst matches = introParagraph.search(\"2016 through the United States\", { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nmatches.items[0].insertText(\"2016, throughout the United States,\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2. Remove the blank paragraph that sits between the two body paragraphs.\nbody.paragraphs.load(\"items\");\nawait context.sync();\nparagraphs = body.paragraphs.items;\nconst blankParagraph = paragraphs[2];\nblankParagraph.delete();\nawait context.sync();\n\n// 3. Replace the \"Add database...\" placeholder paragraph's text with the\n//    real database/table description. Using the paragraph's own range\n//    keeps the \"_GoBack\" bookmark that lives at the end of the paragraph.\nbody.paragraphs.load(\"items\");\nawait context.sync();\nparagraphs = body.paragraphs.items;\nconst placeholderParagraph = paragraphs[2];\n\nconst newParagraphText =\n  \"We are using the Sqlite3 database and have three tables. One table holds the stadium data; this includes: zip codes, city, the year it was built, county and state. Another table holds the census table names and the description of the tables. The remaining tables are created by the program and hold the information obtained by the APIs. Each table holds data obtained from each distinct API used. \";\n\nplaceholderParagraph.getRange().insertText(newParagraphText, Word.InsertLocation.replace);\nawait context.sync();\n\n// 4. Move the \"_GoBack\" bookmark out of that paragraph and into a new,\n//    empty paragraph that follows it.\nbody.paragraphs.load(\"items\");\nawait context.sync();\nparagraphs = body.paragraphs.items;\nconst descriptionParagraph = paragraphs[2];\nconst trailingParagraph = descriptionParagraph.insertParagraph(\"\", Word.InsertLocation.after);\n\ncontext.document.deleteBookmark(\"_GoBack\");\ntrailingParagraph.getRange().insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Expand \"...built from 2011 to 2016 through the United States\" into\n#    \"...built from 2011 to 2016, throughout the United States,\" (adds \", \" + \"out\" + trailing \",\")\n$find = $d.Content\n$find.Find.Execute(\n    \"2016 through the United States\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"2016, throughout the United States,\", 2\n) | Out-Null\n\n# 2. Remove the blank paragraph that sits between the two body paragraphs.\n$d.Paragraphs(3).Range.Delete()\n\n# 3. Replace the \"Add database...\" placeholder paragraph with the real\n#    database/table description.\n$newText = \"We are using the Sqlite3 database and have three tables. One table holds the stadium data; this includes: zip codes, city, the year it was built, county and state. Another table holds the census table names and the description of the tables. The remaining tables are created by the program and hold the information obtained by the APIs. Each table holds data obtained from each distinct API used. \"\n$p3 = $d.Paragraphs(3)\n$p3.Range.Text = $newText\n\n# 4. Move the \"_GoBack\" bookmark out of that paragraph and into a new,\n#    empty paragraph that follows it.\n$endRng = $p3.Range\n$endRng.Collapse(0)\n$endRng.InsertParagraphAfter()\n$d.Bookmarks(\"_GoBack\").Delete()\n$p4 = $d.Paragraphs(4)\n$d.Bookmarks.Add(\"_GoBack\", $p4.Range)\n"}
